# feat(ui): handle empty expense table state
#
# Splits the old single "No expenses stored state" test case (row 7) into
# two distinct test cases:
#   - row 7: the empty-table state itself
#   - row 8 (new): the transition when data is subsequently added
# This is done by inserting a new row at position 8 (shifting the old
# rows 8-35 down to 9-36) and then updating the text content of the two
# "Listing & Rendering" rows (7 & 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$lq = [char]0x201C
$rq = [char]0x201D

# --- Insert a new row at 8, pushing everything below down by one ---
$ws.Rows.Item(8).Insert()

# --- Row 7: "No expenses stored state" (empty table) ---
$ws.Range("A7").Value = "TS_02"
$ws.Range("B7").Value = "TC_02"
$ws.Range("C7").Value = "Listing & Rendering"
$ws.Range("D7").Value = "Sprint 1"
$ws.Range("E7").Value = "UX"
$ws.Range("G7").Value = "1) Clear LocalStorage 2) Reload app"
$ws.Range("H7").Value = "$($lq)No expenses recorded yet$($rq) displayed in table body"
$ws.Range("J7").Value = "Pass"
$ws.Range("I7").Value = "$($lq)No expenses recorded yet$($rq) displayed in table body when table is empty"
$ws.Range("F7").Value = "No expenses stored state when there is no data"
$ws.Rows.Item(7).RowHeight = 72

# --- Row 8 (new): "No expenses stored state" when data is added afterwards ---
$ws.Range("A8").Value = "TS_02"
$ws.Range("B8").Value = "TC_03"
$ws.Range("C8").Value = "Listing & Rendering"
$ws.Range("D8").Value = "Sprint 1"
$ws.Range("E8").Value = "UX"
$ws.Range("F8").Value = "No expenses stored state when data is added"
$ws.Range("G8").Value = "1) Clear LocalStorage 2) Reload app 3) Add some data"
$ws.Range("H8").Value = "$($lq)No expenses recorded yet$($rq) should disappera and new data should be updated"
$ws.Range("I8").Value = "$($lq)No expenses recorded yet$($rq) idsappears in table body and new data is displayed"
$ws.Range("J8").Value = "Pass"
$ws.Rows.Item(8).RowHeight = 72

# --- Append a new trailing filler row (36), matching the style of row 35 ---
$ws.Range("A35:B35").Copy($ws.Range("A36:B36"))

# --- Update the active selection to reflect where the author was working ---
$ws.Range("J8").Select()
